# Insert a new weekly record as row 653, pushing the existing rows
# (653..686) down by one (to 654..687). The new row carries the same
# market / product metadata as its neighbours, but its own date, volume,
# price and origin figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 653; Excel shifts 653:686 down to 654:687
# and copies the formatting (e.g. the date style on column D) from the
# row above, matching the existing sheet layout.
$ws.Rows.Item(653).Insert()

# Populate the newly inserted row 653 with this week's record.
$ws.Range("A653").Value = 10
$ws.Range("B653").Value = "Vega Modelo de Temuco"
$ws.Range("C653").Value = "La Araucanía"
$ws.Range("D653").Value = 45267
$ws.Range("E653").Value = 9
$ws.Range("F653").Value = 100112040
$ws.Range("G653").Value = "Cilantro"
$ws.Range("H653").Value = "Sin especificar"
$ws.Range("I653").Value = "Primera"
$ws.Range("J653").Value = 40
$ws.Range("K653").Value = 7000
$ws.Range("L653").Value = 7000
$ws.Range("M653").Value = 7000
$ws.Range("N653").Value = "$/docena de atados (2 kilos)"
$ws.Range("O653").Value = "Provincia de Cautín"
$ws.Range("P653").Value = 3500
$ws.Range("Q653").Value = 2
$ws.Range("R653").Value = "Hortaliza"
